$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Fgf1","Cspg4","ECs",2,0.6666666666666666,0.265349,0.7960469999999999,0.1498685997319469,0.1498685997319469,3,1,1.415809333333333,4.247427999999999,0.1174285988148896,0.1174285988148896,0.3756835907906666,3.381152317115999,0.01759885967287206,0.01759885967287207),
    @("ECs","Fgf1","Cspg4","FAPs",2,0.6666666666666666,0.265349,0.7960469999999999,0.1498685997319469,0.1498685997319469,3,1,6.180051333333334,18.540154,0.5125794495003263,0.5125794495003263,1.639870441248667,14.758833971238,0.07681956434798606,0.07681956434798606),
    @("ECs","Fgf1","Cspg4","Inflammatory-Mac",2,0.6666666666666666,0.265349,0.7960469999999999,0.1498685997319469,0.1498685997319469,1,0.3333333333333333,0.130008,0.390024,0.01078298957020072,0.01078298957020072,0.034497492792,0.310477435128,0.00161603154781017,0.00161603154781017),
    @("ECs","Fgf1","Cspg4","MuSCs",2,0.6666666666666666,0.265349,0.7960469999999999,0.1498685997319469,0.1498685997319469,3,1,3.614693333333333,10.84408,0.2998061697188437,0.2998061697188436,0.9591552613066666,8.63239735176,0.04493153084676151,0.0449315308467615),
    @("ECs","Fgf1","Cspg4","Resolving-Mac",2,0.6666666666666666,0.265349,0.7960469999999999,0.1498685997319469,0.1498685997319469,3,1,0.7162056666666667,2.148617,0.05940279239573968,0.05940279239573967,0.1900444574443333,1.710400116999,0.008902613316517047,0.008902613316517045),
    @("FAPs","Fgf1","Cspg4","ECs",3,1,1.003400666666667,3.010202,0.5667187473230926,0.5667187473230925,3,1,1.415809333333333,4.247427999999999,0.1174285988148896,0.1174285988148896,1.420624028939555,12.785616260456,0.06654898842028024,0.06654898842028023),
    @("FAPs","Fgf1","Cspg4","FAPs",3,1,1.003400666666667,3.010202,0.5667187473230926,0.5667187473230925,3,1,6.180051333333334,18.540154,0.5125794495003263,0.5125794495003263,6.20106762790089,55.80960865110799,0.2904883835243853,0.2904883835243853),
    @("FAPs","Fgf1","Cspg4","Inflammatory-Mac",3,1,1.003400666666667,3.010202,0.5667187473230926,0.5667187473230925,1,0.3333333333333333,0.130008,0.390024,0.01078298957020072,0.01078298957020072,0.130450113872,1.174051024848,0.006110922341622126,0.006110922341622125),
    @("FAPs","Fgf1","Cspg4","MuSCs",3,1,1.003400666666667,3.010202,0.5667187473230926,0.5667187473230925,3,1,3.614693333333333,10.84408,0.2998061697188437,0.2998061697188436,3.626985700462222,32.64287130416,0.1699057769427976,0.1699057769427975),
    @("FAPs","Fgf1","Cspg4","Resolving-Mac",3,1,1.003400666666667,3.010202,0.5667187473230926,0.5667187473230925,3,1,0.7162056666666667,2.148617,0.05940279239573968,0.05940279239573967,0.7186412434037778,6.467771190634,0.03366467609400732,0.03366467609400731),
    @("MuSCs","Fgf1","Cspg4","ECs",3,1,0.4790693333333333,1.437208,0.2705774288246194,0.2705774288246194,3,1,1.415809333333333,4.247427999999999,0.1174285988148896,0.1174285988148896,0.6782708334471109,6.104437501023998,0.03177352833781059,0.03177352833781059),
    @("MuSCs","Fgf1","Cspg4","FAPs",3,1,0.4790693333333333,1.437208,0.2705774288246194,0.2705774288246194,3,1,6.180051333333334,18.540154,0.5125794495003263,0.5125794495003263,2.960673072225778,26.646057650032,0.1386924295141372,0.1386924295141372),
    @("MuSCs","Fgf1","Cspg4","Inflammatory-Mac",3,1,0.4790693333333333,1.437208,0.2705774288246194,0.2705774288246194,1,0.3333333333333333,0.130008,0.390024,0.01078298957020072,0.01078298957020072,0.06228284588799999,0.5605456129919999,0.0029176335929476,0.0029176335929476),
    @("MuSCs","Fgf1","Cspg4","MuSCs",3,1,0.4790693333333333,1.437208,0.2705774288246194,0.2705774288246194,3,1,3.614693333333333,10.84408,0.2998061697188437,0.2998061697188436,1.731688725404444,15.58519852864,0.0811207825482822,0.08112078254828219),
    @("MuSCs","Fgf1","Cspg4","Resolving-Mac",3,1,0.4790693333333333,1.437208,0.2705774288246194,0.2705774288246194,3,1,0.7162056666666667,2.148617,0.05940279239573968,0.05940279239573967,0.3431121712595556,3.088009541336,0.0160730548314419,0.0160730548314419),
    @("Resolving-Mac","Fgf1","Cspg4","ECs",1,0.3333333333333333,0.02272533333333333,0.068176,0.01283522412034115,0.01283522412034115,3,1,1.415809333333333,4.247427999999999,0.1174285988148896,0.1174285988148896,0.03217473903644444,0.289572651328,0.001507222383926735,0.001507222383926735),
    @("Resolving-Mac","Fgf1","Cspg4","FAPs",1,0.3333333333333333,0.02272533333333333,0.068176,0.01283522412034115,0.01283522412034115,3,1,6.180051333333334,18.540154,0.5125794495003263,0.5125794495003263,0.1404437265671111,1.263993539104,0.006579072113817775,0.006579072113817775),
    @("Resolving-Mac","Fgf1","Cspg4","Inflammatory-Mac",1,0.3333333333333333,0.02272533333333333,0.068176,0.01283522412034115,0.01283522412034115,1,0.3333333333333333,0.130008,0.390024,0.01078298957020072,0.01078298957020072,0.002954475136,0.026590276224,0.0001384020878208273,0.0001384020878208273),
    @("Resolving-Mac","Fgf1","Cspg4","MuSCs",1,0.3333333333333333,0.02272533333333333,0.068176,0.01283522412034115,0.01283522412034115,3,1,3.614693333333333,10.84408,0.2998061697188437,0.2998061697188436,0.08214511089777778,0.73930599808,0.003848079381002394,0.003848079381002393),
    @("Resolving-Mac","Fgf1","Cspg4","Resolving-Mac",1,0.3333333333333333,0.02272533333333333,0.068176,0.01283522412034115,0.01283522412034115,3,1,0.7162056666666667,2.148617,0.05940279239573968,0.05940279239573967,0.01627601251022223,0.146484112592,0.0007624481537734155,0.0007624481537734154)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}
